$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$dateFmt = "[$-409]d/mmm/yyyy;@"
$rupee = [char]0x20B9
$amtFmt = '"' + $rupee + '"#,##0;"' + $rupee + '"\-#,##0'
$newDate = 44795

$rows = @(
  @{ Row=295; B="KA51MN2772"; C="DUSTER";      D="SUSPENSION";              E="WORK DONE DELIVERED"; F=38000; G="P PAY" },
  @{ Row=296; B="MP50BC8265"; C="INNOVA";      D="HEAD LIGHT BULB CHANGE"; E="WORK DONE DELIVERED"; F=1020;  G="P PAY" },
  @{ Row=297; B="KA50P8555";  C="I10";         D="BODY SHOP";               E="WORK DONE";           F=11308; G="  INSURANCE" },
  @{ Row=298; B="KA03AD4806"; C="VERITO";      D="BRAKE DISC CHANGE";       E="WORK DONE DELIVERED"; F=4500;  G="P PAY" },
  @{ Row=299; B="KA03MS5800"; C="I10 ";        D="PMS  ";                   E="WORK IN PROGRESS" },
  @{ Row=300; B="KA03MT2662"; C="POLO";        D="PIPE CHANGE";             E="WORK IN PROGRESS" },
  @{ Row=301; B="KA25P8050";  C="NANO";        D="PMS";                     E="WORK IN PROGRESS" },
  @{ Row=302; B="KA03MJ4740"; C="HONDA CITY";  D="PMS                                      WW"; E="WORK DONE DELIVERED"; F=3233;  G="CREDIT" },
  @{ Row=303; B="KA03NA8797 ";C="JAZZ";        D="PMS";                     E="WORK DONE DELIVERED"; F=11389 },
  @{ Row=304; B="KA03NC9110"; C="TIAGO";       D="PMS";                     E="WORK DONE DELIVERED"; F=4494;  G="G PAY" }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $newDate
    $ws.Range("A$row").NumberFormat = $dateFmt
    $ws.Range("B$row").Value = $r.B
    $ws.Range("C$row").Value = $r.C
    $ws.Range("D$row").Value = $r.D
    $ws.Range("E$row").Value = $r.E
    if ($r.ContainsKey("F")) {
        $ws.Range("F$row").Value = $r.F
        $ws.Range("F$row").NumberFormat = $amtFmt
    }
    if ($r.ContainsKey("G")) {
        $ws.Range("G$row").Value = $r.G
    }
}

$ws.Range("A1:G304").Select()
$ws.Application.ActiveWindow.ScrollRow = 280
$ws.Range("G302").Select()
